$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D26").Value = "생성 모델의 새로운 흐름 확산 모델(Diffusion model)에 관하여"

$ws.Range("D37").Value = "[Paper Review] Improving Language Models by Retrieving from Trillions of Tokens"

$ws.Range("D50").Value = "Mahalanobis distance"
$ws.Range("E50").Value = "http://incredible.egloos.com/7539231"

$ws.Range("D51").Value = "[css] 요소 사이에 구분선 넣고 싶다면? ex) 게시판 | 1:1문의 | 회원가입"
$ws.Range("E51").Value = "https://bskyvision.com/1260"
